$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates ---
$rushing = $wb.Worksheets.Item("Rushing")

# J.Fromm
$rushing.Range("C2").Value = 3

# D.Booker
$rushing.Range("C4").Value = 56
$rushing.Range("D4").Value = 28
$rushing.Range("E4").Value = 13

# G.Brightwell
$rushing.Range("C5").Value = 66
$rushing.Range("D5").Value = 56
$rushing.Range("E5").Value = 15

# --- Receiving sheet updates ---
$receiving = $wb.Worksheets.Item("Receiving")

# D.Booker
$receiving.Range("C2").Value = 33

# G.Brightwell
$receiving.Range("C3").Value = 41
$receiving.Range("D3").Value = 37

# K.Golladay
$receiving.Range("C6").Value = 48

# P.Cooper
$receiving.Range("C13").Value = 11

# D.Pettis
$receiving.Range("C14").Value = 2
$receiving.Range("D14").Value = 2

# D.Sills
$receiving.Range("C15").Value = 64
$receiving.Range("D15").Value = 43

# K.Rudolph
$receiving.Range("C16").Value = 29
